# Applies cryptos list price/volume updates per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$sub3 = [char]0x2083

$cellD = $ws.Range("D2")
$cellD.NumberFormat = "@"
$cellD.Value = "71.038.92"
$ws.Range("E2").Value = "  -2.29%  "

$cellD = $ws.Range("D3")
$cellD.NumberFormat = "@"
$cellD.Value = "3.855.30"
$ws.Range("E3").Value = "  -2.86%  "

$cellD = $ws.Range("D4")
$cellD.NumberFormat = "@"
$cellD.Value = "0.999"
$ws.Range("E4").Value = "  -0.23%  "

$cellD = $ws.Range("D5")
$cellD.NumberFormat = "@"
$cellD.Value = "591.58"
$ws.Range("E5").Value = "  +0.98%  "

$cellD = $ws.Range("D6")
$cellD.NumberFormat = "@"
$cellD.Value = "166.71"
$ws.Range("E6").Value = "  +5.79%  "

$cellD = $ws.Range("D7")
$cellD.NumberFormat = "@"
$cellD.Value = "0.670"
$ws.Range("E7").Value = "  -1.45%  "

$cellD = $ws.Range("D8")
$cellD.NumberFormat = "@"
$cellD.Value = "0.999"
$ws.Range("E8").Value = "  -0.02%  "

$cellD = $ws.Range("D9")
$cellD.NumberFormat = "@"
$cellD.Value = "0.749"
$ws.Range("E9").Value = "  +0.25%  "

$cellD = $ws.Range("D10")
$cellD.NumberFormat = "@"
$cellD.Value = "0.173"
$ws.Range("E10").Value = "  +3.65%  "

$cellD = $ws.Range("D11")
$cellD.NumberFormat = "@"
$cellD.Value = "53.21"
$ws.Range("E11").Value = "  -1.35%  "

$cellD = $ws.Range("D12")
$cellD.NumberFormat = "@"
$cellD.Value = "0.0000320"
$ws.Range("E12").Value = "  +0.92%  "

$cellD = $ws.Range("D13")
$cellD.NumberFormat = "@"
$cellD.Value = "11.13"
$ws.Range("E13").Value = "  +2.71%  "

$cellD = $ws.Range("D14")
$cellD.NumberFormat = "@"
$cellD.Value = "4.472.92"
$ws.Range("E14").Value = "  -2.95%  "

$cellD = $ws.Range("D15")
$cellD.NumberFormat = "@"
$cellD.Value = "3.883.93"
$ws.Range("E15").Value = "  -2.11%  "

$cellD = $ws.Range("D16")
$cellD.NumberFormat = "@"
$cellD.Value = "20.65"
$ws.Range("E16").Value = "  +1.21%  "

$cellD = $ws.Range("D17")
$cellD.NumberFormat = "@"
$cellD.Value = "13.75"
$ws.Range("E17").Value = "  -1.63%  "

$ws.Range("E18").Value = "  -5.87%  "

$cellD = $ws.Range("D20")
$cellD.NumberFormat = "@"
$cellD.Value = "70.832.23"
$ws.Range("E20").Value = "  -2.36%  "

$cellD = $ws.Range("D21")
$cellD.NumberFormat = "@"
$cellD.Value = "432.14"
$ws.Range("E21").Value = "  +0.09%  "

$cellD = $ws.Range("D22")
$cellD.NumberFormat = "@"
$cellD.Value = "4.71"
$ws.Range("E22").Value = "  +0.38%  "

$cellD = $ws.Range("D23")
$cellD.NumberFormat = "@"
$cellD.Value = "94.04"
$ws.Range("E23").Value = "  -1.91%  "

$cellD = $ws.Range("D24")
$cellD.NumberFormat = "@"
$cellD.Value = "3.26"
$ws.Range("E24").Value = "  -4.86%  "

$cellD = $ws.Range("D25")
$cellD.NumberFormat = "@"
$cellD.Value = "13.70"
$ws.Range("E25").Value = "  -4.21%  "

$cellD = $ws.Range("D26")
$cellD.NumberFormat = "@"
$cellD.Value = "4.12"
$ws.Range("E26").Value = "  -6.92%  "

$cellD = $ws.Range("D27")
$cellD.NumberFormat = "@"
$cellD.Value = "10.82"
$ws.Range("E27").Value = "  -3.31%  "

$cellD = $ws.Range("D28")
$cellD.NumberFormat = "@"
$cellD.Value = "5.91"
$ws.Range("E28").Value = "  -0.40%  "

$cellD = $ws.Range("D29")
$cellD.NumberFormat = "@"
$cellD.Value = "10.14"
$ws.Range("E29").Value = "  -5.67%  "

$cellD = $ws.Range("D30")
$cellD.NumberFormat = "@"
$cellD.Value = "34.87"
$ws.Range("E30").Value = "  -3.96%  "

$cellD = $ws.Range("D31")
$cellD.NumberFormat = "@"
$cellD.Value = "7.87"
$ws.Range("E31").Value = "  +0.40%  "

$cellD = $ws.Range("D32")
$cellD.NumberFormat = "@"
$cellD.Value = "49.75"
$ws.Range("E32").Value = "  -0.77%  "

$cellD = $ws.Range("D33")
$cellD.NumberFormat = "@"
$cellD.Value = "13.44"
$ws.Range("E33").Value = "  -0.99%  "

$ws.Range("E34").Value = "  -5.15%  "

$cellD = $ws.Range("D35")
$cellD.NumberFormat = "@"
$cellD.Value = "68.69"
$ws.Range("E35").Value = "  +0.15%  "

$cellD = $ws.Range("D36")
$cellD.NumberFormat = "@"
$cellD.Value = "0.0${sub3}0973"
$ws.Range("E36").Value = "  +13.51%  "

$cellD = $ws.Range("D37")
$cellD.NumberFormat = "@"
$cellD.Value = "617.00"
$ws.Range("E37").Value = "  -9.08%  "

$cellD = $ws.Range("D38")
$cellD.NumberFormat = "@"
$cellD.Value = "0.417"
$ws.Range("E38").Value = "  -4.41%  "

$ws.Range("E39").Value = "  +0.18%  "

$ws.Range("E40").Value = "  -0.23%  "

$ws.Range("E41").Value = "  -2.22%  "

$cellD = $ws.Range("D42")
$cellD.NumberFormat = "@"
$cellD.Value = "0.141"
$ws.Range("E42").Value = "  -3.03%  "

$cellD = $ws.Range("D43")
$cellD.NumberFormat = "@"
$cellD.Value = "3.23"
$ws.Range("E43").Value = "  +35.67%  "

$cellD = $ws.Range("D44")
$cellD.NumberFormat = "@"
$cellD.Value = "0.0464"
$ws.Range("E44").Value = "  -4.38%  "

$cellD = $ws.Range("D45")
$cellD.NumberFormat = "@"
$cellD.Value = "10.05"
$ws.Range("E45").Value = "  -8.15%  "

$cellD = $ws.Range("D46")
$cellD.NumberFormat = "@"
$cellD.Value = "2.63"
$ws.Range("E46").Value = "  -1.87%  "

$ws.Range("E47").Value = "  -3.71%  "

$cellD = $ws.Range("D48")
$cellD.NumberFormat = "@"
$cellD.Value = "3.38"
$ws.Range("E48").Value = "  -0.43%  "

$cellD = $ws.Range("D49")
$cellD.NumberFormat = "@"
$cellD.Value = "2.818.77"
$ws.Range("E49").Value = "  +2.14%  "

$ws.Range("E50").Value = "  -18.47%  "

$cellD = $ws.Range("D51")
$cellD.NumberFormat = "@"
$cellD.Value = "0.000272"
$ws.Range("E51").Value = "  +1.23%  "
